$d = $word.ActiveDocument

function Insert-EmptyParaAfter($afterIndex) {
    $rng = $d.Paragraphs.Item($afterIndex).Range
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    return ($afterIndex + 1)
}

function Insert-TextParaAfter($afterIndex, $text) {
    $newIndex = Insert-EmptyParaAfter $afterIndex
    $p = $d.Paragraphs.Item($newIndex)
    $p.Range.Text = $text
    return $newIndex
}

function Set-ParaBold($index) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Font.Bold = 1
}

# Hunk 1: insert 15 empty paragraphs before the "Return View(), RedirectToAction()..." heading
# (originally paragraph 41 is the empty paragraph right before that heading, which was originally paragraph 42)
$idx = 41
for ($i = 0; $i -lt 15; $i++) {
    $idx = Insert-EmptyParaAfter $idx
}

# Hunk 2: insert new Database/EF Core documentation section after the paragraph
# "return RedirectToAction(...)" which was originally paragraph 51, but has shifted
# to index 66 (51 + 15) because of the 15 paragraphs inserted for Hunk 1 above.
$idx = 66
$boldIndexes = New-Object System.Collections.ArrayList

$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "Database"
[void]$boldIndexes.Add($idx)
$idx = Insert-TextParaAfter $idx "Biz səhifədə göstərməli olduğumuz dataları Database daxilində saxlayırıq. Database’lər isə bizim məqsədimizə və ya istəyimizə uyğun olaraq dəyişir."
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "DbContext"
[void]$boldIndexes.Add($idx)
$idx = Insert-TextParaAfter $idx "İlk olaraq DbContext class’dan miras alan AppDbContext class’ı yaratmalıyıq. DbContext’dən miras ala bilmək üçün NuGet Packages’dən Microsoft.EntityFrameworkCore package yüklənməlidir. AppDbContext class daxilində Constructor yaradılmalı və parametr olaraq DbContextOptions<AppDbContext> type’dan olan options object’i qəbul etməlidir daha sonra  base class’a(DbContext class’a) options object’i göndərməlidir. "
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "DbSet<T>"
[void]$boldIndexes.Add($idx)
$idx = Insert-TextParaAfter $idx "Beləliklə biz Database daxilində yaranmasını istədiyimiz Table’lar üçün AppDbContext class daxilində generic olan public DbSet<EntityName> EntityNames {get;set;} property’ləri yaratmalıyıq. Yaratdığımız hər bir DbSet property Database daxilində bir table’a qarşılıq gəlir bu səbəbdən generic type olaraq model class’larımızı yazırıq."
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "AddDbContext"
[void]$boldIndexes.Add($idx)
$idx = Insert-TextParaAfter $idx "Biz Program.cs daxilində Services hissədə Dependency Injection baş verə bilsin və Database ilə əlaqə yarana bilsin deyə bildirməliyik ki DbContext’dən istifadə etmişik. builder.Services.AddDbContext<AppDbContext>(opt=>{ opt.UseSqlServer(connectionString) });"
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "UseSqlServer"
[void]$boldIndexes.Add($idx)
$idx = Insert-TextParaAfter $idx "UseSqlServer’dən istifadə edə bilməyimiz üçün"
$idx = Insert-TextParaAfter $idx "NuGet Packages’dən MsSql üçün istifadə olunan Microsoft.EntityFrameworkCore.SqlServer package yüklənməlidir."
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "Connection string"
[void]$boldIndexes.Add($idx)
$idx = Insert-TextParaAfter $idx "Connection string məxfi data olduğu üçün onu kənardan çağırmalıyıq. appsettings.json file daxilində json şəklində saxlayacayıq. Connection string aşağıdaki şəkildə olmalıdır."
$idx = Insert-TextParaAfter $idx "“ConnectionStrings” :{“Default” : “Server:ServerName;Database=DatabaseName;Trustesd_Connection=True;”}"
$idx = Insert-TextParaAfter $idx ".cs file daxilində bir neçə yolla çağrıla bilər. Bunlardan biri Configuration class’dan istifadə etməkdir. builder.Configuration[“Key:Value”] məntiqi ilə çağırıla bilər. builder.GetConnectionString(“Value”); yazaraq da eyni işi görə bilərik."
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "Migrations"
[void]$boldIndexes.Add($idx)
$idx = Insert-TextParaAfter $idx "Biz code hissədə etdiyimiz dəyişikliklərin Database’də öz əksini tapa bilməsi üçün terminal olaraq həm Package Manager Console’dan həm də Developer PowerShell’dən istifadə edə bilərik. Lakin NuGet Packages vasitəsi ilə Microsoft.EntityFrameworkCore.Tools package yüklənməlidir. Package Manager Console’da migration’ları yaradan zaman Default Project olaraq AppDbContext’in yerləşdiyi layer’ı seçirik. Migration’ların yarana bilməsi üçün əlavə olaraq Microsoft.EntityFrameworkCore.Design package yüklənməlidir."
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "add-migration “MigrationName” yazaraq edilən dəyişiklikləri yeni yaratdığımız migration daxilinə əlavə edirik."
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "update-database vasitəsi ilə biz yaratdığımız yeni migration’u database’ə əlavə edirik."
$idx = Insert-TextParaAfter $idx "update-database MigrationName vasitəsi ilə biz database’i daha əvvəldən yaratdığımız migration halına gətirə bilərik."
$idx = Insert-EmptyParaAfter $idx
$idx = Insert-TextParaAfter $idx "remove-migration vasitəsi ilə update olunmamış migration’u silirik. Update olunan migration’u silmək istədiyimiz halda isə ondan əvvəlki migration’a update edərək silmək istədiyimiz migration’u  remove-migration yazaraq silirik. Manual olaraq gedib migration folder daxilindən migration’u silmək məsləhət görülmür."
$idx = Insert-EmptyParaAfter $idx

foreach ($bi in $boldIndexes) {
    Set-ParaBold $bi
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
